$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 10:50:50.582329",
    "2021-10-05 10:50:50.582340",
    "2021-10-05 10:50:50.582343",
    "2021-10-05 10:50:50.582346",
    "2021-10-05 10:50:50.582349",
    "2021-10-05 10:50:50.582351",
    "2021-10-05 10:50:50.582354",
    "2021-10-05 10:50:50.582357",
    "2021-10-05 10:50:50.582359",
    "2021-10-05 10:50:50.582362",
    "2021-10-05 10:50:50.582364",
    "2021-10-05 10:50:50.582367",
    "2021-10-05 10:50:50.582369",
    "2021-10-05 10:50:50.582372",
    "2021-10-05 10:50:50.582375",
    "2021-10-05 10:50:50.582377",
    "2021-10-05 10:50:50.582380",
    "2021-10-05 10:50:50.582383",
    "2021-10-05 10:50:50.582385",
    "2021-10-05 10:50:50.582388"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
